$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended to the bottom of the Adafruit IO feed data (row 17),
# matching the existing schema: Timestamp, Feed Key, Value, Latitude, Longitude, Elevation
$row = 17

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# The "Value" column holds numeric-looking readings but the sheet stores them
# as text (see existing rows), so force text formatting before assigning the
# value to avoid Excel auto-converting "25" into the number 25. Clear the
# formatting afterwards so the cell doesn't pick up a stray explicit style
# and stays consistent with the rest of the (unstyled) sheet.
$valueCell = $ws.Cells.Item($row, 3)
$valueCell.NumberFormat = "@"
$valueCell.Value = "25"
$valueCell.ClearFormats()

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
